$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column F (dSF) keyed by row number, per the diff.
# Rows not listed below (8, 9, 11, 13, 34, 36) are unchanged.
$updates = @{
    2  = -5
    3  = -2
    4  = -3
    5  = 2
    6  = 1
    7  = 0
    10 = -2
    12 = -3
    14 = 3
    15 = -3
    16 = 2
    17 = 6
    18 = -1
    19 = -2
    20 = -1
    21 = -2
    22 = 4
    23 = 2
    24 = -3
    25 = -2
    26 = 1
    27 = -2
    28 = -1
    29 = -2
    30 = -2
    31 = 1
    32 = 1
    33 = -1
    35 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
